$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert five new "Body Text" (TextBody) paragraphs right after the
#    paragraph ending "...Bet breakpoint to midwinter or to change to
#    ectodormancy?" and before the (pre-existing) empty paragraph that leads
#    into the "Physiology/ mechanisms" heading.
# ---------------------------------------------------------------------------

$anchorText = "Could I just use a break point model to get a slope for autumn and a slope for spring? Bet breakpoint to midwinter or to change to ectodormancy?"

$find = $d.Content
$found = $find.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertionPoint = $find
$insertionPoint.Collapse(0)

for ($i = 0; $i -lt 5; $i++) {
    $insertionPoint.InsertParagraphAfter()
}

# Re-locate the anchor paragraph so we can address the five freshly minted
# paragraphs that now follow it (they inherited the "TextBody" style).
$find2 = $d.Content
$null = $find2.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorParaIndex = $find2.Paragraphs(1).Index

$newPara1 = $d.Paragraphs($anchorParaIndex + 1)
$newPara2 = $d.Paragraphs($anchorParaIndex + 2)
$newPara3 = $d.Paragraphs($anchorParaIndex + 3)
$newPara4 = $d.Paragraphs($anchorParaIndex + 4)
$newPara5 = $d.Paragraphs($anchorParaIndex + 5)

# Paragraph 1: stays empty.

# Paragraph 2.
$newPara2.Range.InsertBefore("Maybe I could have a transformation so that teh relationship is not linear? Vines acheive winter hardiness using supercooling of intracellular water (see Kovaleski et al 2018 intro), and the maximum hardiness possible with this mechanism is 40 degrees C. (Biggs, 1953 in Kovaleski et al 2018 ). Maybe I could use this number in teh non linear transformation?  ")

# Paragraph 3 (originally three runs; text concatenated here).
$newPara3.Range.InsertBefore("Maybe uses a logistic regression? Kovaleski et al 2018 used a logistic regression to get deaclimation rate. Can it be used for non binary data?")

# Paragraph 4: stays empty.

# Paragraph 5 (originally two runs; text concatenated here).
$newPara5.Range.InsertBefore("Kovaleski et al 2018 found that different varieties of winegrape had different rates of deacclimation ")

# ---------------------------------------------------------------------------
# 2. Remove the leading, standalone space-only run immediately before
#    "Energy is necessary to drive acclimation..." in the Physiology intro
#    paragraph.
# ---------------------------------------------------------------------------

$energyAnchor = "Energy is necessary to drive acclimation"
$find3 = $d.Content
$null = $find3.Find.Execute($energyAnchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$energyPara = $find3.Paragraphs(1).Range
$leadingChar = $d.Range($energyPara.Start, $energyPara.Start + 1)
if ($leadingChar.Text -eq " ") {
    $leadingChar.Delete()
}

# ---------------------------------------------------------------------------
# 3. Flip the "Normal" style's overflowPunct (HangingPunctuation) compat
#    setting from true to false.
# ---------------------------------------------------------------------------

$normalStyle = $d.Styles("Normal")
$normalStyle.ParagraphFormat.HangingPunctuation = $false
